$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - AREPD
$ws.Range("B2").Value = 5.201
$ws.Range("C2").Value = 0.443
$ws.Range("D2").Value = 91.48999999999999
$ws.Range("E2").Value = 27.9697

# Row 3 - AV-MCPS
$ws.Range("B3").Value = 1.219
$ws.Range("C3").Value = 0.421
$ws.Range("D3").Value = 65.44
$ws.Range("E3").Value = 19.8846

# Row 4 - Block Bootstrapping
$ws.Range("B4").Value = 6.212
$ws.Range("C4").Value = 0.457
$ws.Range("D4").Value = 92.65000000000001
$ws.Range("E4").Value = 29.5439

# Row 5 - DeepAR
$ws.Range("B5").Value = 0.829
$ws.Range("C5").Value = 0.401
$ws.Range("D5").Value = 51.63
$ws.Range("E5").Value = 14.0055

# Row 6 - EnCQR-LSTM
$ws.Range("B6").Value = 3.838
$ws.Range("C6").Value = 0.706
$ws.Range("D6").Value = 81.61
$ws.Range("E6").Value = 31.8311

# Row 7 - LSPM
$ws.Range("B7").Value = 0.656
$ws.Range("C7").Value = 0.425
$ws.Range("D7").Value = 35.23
$ws.Range("E7").Value = 19.4984

# Row 8 - LSPMW
$ws.Range("B8").Value = 1.621
$ws.Range("C8").Value = 0.476
$ws.Range("D8").Value = 70.66
$ws.Range("E8").Value = 22.9849

# Row 9 - MCPS
$ws.Range("B9").Value = 1.184
$ws.Range("C9").Value = 0.422
$ws.Range("D9").Value = 64.36
$ws.Range("E9").Value = 18.5696

# Row 10 - Sieve Bootstrap
$ws.Range("B10").Value = 0.371
$ws.Range("C10").Value = 0.377
$ws.Range("D10").Value = -1.56
$ws.Range("E10").Value = 0.7361
$ws.Range("F10").Value = 0.4616500615185113
